$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Price column values are forced to Text format so that numeric-looking
# strings (e.g. "250.71") are not converted into floating point numbers,
# matching the original inline-string content of the cells. The cell
# style is restored to "Normal" afterwards so no extra formatting is introduced.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "42.623.54"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.267.04"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.06%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "250.71"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.634"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "75.72"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +6.26%  "
$ws.Range("E8").Value = "  -0.02%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.644"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -2.07%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "40.27"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +4.32%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0971"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "7.32"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.92%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.106"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.20%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.607.93"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.68%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "15.02"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.866"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.269.55"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "42.533.54"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.46%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0994"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("E20").Value = "  -2.18%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "72.05"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.37%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "234.27"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.48%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "2.15"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("E25").Value = "  +0.09%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "11.24"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  -1.96%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.12"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.75%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "167.87"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "21.11"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.22%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "6.47"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.71%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.0855"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +7.03%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "32.16"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.12%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.124"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("E35").Value = "  +1.28%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "4.53"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.67%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "4.73"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  -5.11%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "13.58"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +6.33%  "
$ws.Range("E40").Value = "  -2.01%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.85"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("E42").Value = "  +1.79%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "61.48"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.35%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "8.87"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -3.45%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "106.20"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +11.82%  "
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("E47").Value = "  -4.03%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("E51").Value = "  -2.38%  "
